$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the divisor in the D-column formulas from /10 to /15 for rows 32, 33, 37, 41, 45
$rows = @(32, 33, 37, 41, 45)
foreach ($r in $rows) {
    $ws.Range("D$r").Formula = "=SUM(`$D`$52:`$D`$85)*(35/34) / 15"
}


# Update the view: scroll so A28 is the top-left visible cell, and select M17
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M17").Select()
